$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update FirstName column (A2:A4)
$ws.Range("A2").Value = "Adeyy"
$ws.Range("A3").Value = "Poef"
$ws.Range("A4").Value = "Lofd"

# Update Username column (E2:E4)
$ws.Range("E2").Value = "crewr354"
$ws.Range("E3").Value = "ftry3"
$ws.Range("E4").Value = "suhn35"

# Move selection to C2
$ws.Range("C2").Select() | Out-Null
